$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Execution-Times")

# Updated Spark Execution Time (seconds) results for the iterative method.
$values = @{
    2  = 16.8
    3  = 1.19
    4  = 3.08
    5  = 0.93
    6  = 0.96
    7  = 15.88
    8  = 1.49
    9  = 3.71
    10 = 0.95
    11 = 0.9
    12 = 20.79
    13 = 1.33
    14 = 0.82
    15 = 1.08
    16 = 0.8
    17 = 15.31
    18 = 1.31
    19 = 0.9
    20 = 1.13
    21 = 1.02
    22 = 19.82
    23 = 1.65
    24 = 0.91
    25 = 0.75
    26 = 0.79
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# A handful of rows also switch from the "General" number style to the
# "0.00" (2 decimal place) number style used elsewhere in the column.
$restyledRows = @(2, 11, 16, 19)
foreach ($row in $restyledRows) {
    $ws.Cells.Item($row, 4).NumberFormat = "0.00"
}

# Reflect the zoom/scroll/selection state left behind on the
# "Execution-Times" sheet, then restore the "Average-Execution-Time"
# sheet as the active tab (matching the workbook's saved view state).
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 70
$ws.Range("AC22").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("Average-Execution-Time")
$ws2.Select() | Out-Null
